$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of mod-count data for 2025/11/16.
$newRow = 7

# Write the date as literal text (matching the existing "YYYY/MM/DD" text
# cells above it) instead of letting Excel auto-convert it to a date serial.
$ws.Cells.Item($newRow, 1).NumberFormat = "@"
$ws.Cells.Item($newRow, 1).Value = "2025/11/16"
$ws.Cells.Item($newRow, 1).ClearFormats()

$ws.Cells.Item($newRow, 2).Value = "逃离鸭科夫"
$ws.Cells.Item($newRow, 3).Value = 1139

# Match the formatting (centered alignment) used by the preceding data rows.
$ws.Range("A6:C6").Copy()
$ws.Range("A7:C7").PasteSpecial(-4122)
$excel.CutCopyMode = $false
